# Adds two new columns, I ("I0") and J ("IF"), to the existing table in row 1..60.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new labels in I1/J1, matching H1 formatting (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-60: numeric I/J values, unstyled like the existing C:H data columns.
$ijData = @{
    2 = @(6, 6)
    3 = @(6, 7)
    4 = @(8, 8)
    5 = @(9, 9)
    6 = @(10, 10)
    7 = @(7, 7)
    8 = @(10, 10)
    9 = @(8, 8)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(7, 7)
    17 = @(9, 9)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(6, 7)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(9, 9)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(8, 8)
    32 = @(9, 9)
    33 = @(8, 8)
    34 = @(6, 6)
    35 = @(8, 8)
    36 = @(8, 8)
    37 = @(7, 7)
    38 = @(9, 9)
    39 = @(9, 9)
    40 = @(6, 7)
    41 = @(7, 7)
    42 = @(8, 8)
    43 = @(8, 8)
    44 = @(5, 6)
    45 = @(7, 7)
    46 = @(6, 6)
    47 = @(6, 6)
    48 = @(7, 7)
    49 = @(2, 2)
    50 = @(6, 6)
    51 = @(7, 7)
    52 = @(5, 5)
    53 = @(7, 7)
    54 = @(6, 6)
    55 = @(9, 9)
    56 = @(7, 7)
    57 = @(6, 6)
    58 = @(8, 8)
    59 = @(7, 7)
    60 = @(6, 6)
}

foreach ($r in $ijData.Keys) {
    $vals = $ijData[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
